$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to Text so numeric-looking strings
# (e.g. "0.9997", "23.390.87") are stored verbatim and not coerced to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '23.390.87'
$ws.Range('E2').Value = '  -0.20%  '
$ws.Range('D3').Value = '1.630.43'
$ws.Range('E3').Value = '  -0.71%  '
$ws.Range('D4').Value = '0.9997'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '0.9994'
$ws.Range('E5').Value = '  -0.03%  '
$ws.Range('D6').Value = '302.23'
$ws.Range('E6').Value = '  -0.96%  '
$ws.Range('D7').Value = '0.3768'
$ws.Range('E7').Value = '  +0.90%  '
$ws.Range('D8').Value = '51.83'
$ws.Range('E8').Value = '  -0.97%  '
$ws.Range('E9').Value = '  -0.12%  '
$ws.Range('D10').Value = '0.08174'
$ws.Range('E10').Value = '  +0.60%  '
$ws.Range('E11').Value = '  -2.73%  '
$ws.Range('D12').Value = '0.9996'
$ws.Range('E12').Value = '  -0.03%  '
$ws.Range('D13').Value = '22.24'
$ws.Range('E13').Value = '  -2.59%  '
$ws.Range('E14').Value = '  -2.08%  '
$ws.Range('D15').Value = '7.319'
$ws.Range('E15').Value = '  +0.39%  '
$ws.Range('D16').Value = '0.00001240'
$ws.Range('D17').Value = '1.621.54'
$ws.Range('E17').Value = '  -0.58%  '
$ws.Range('D18').Value = '94.74'
$ws.Range('E18').Value = '  +0.32%  '
$ws.Range('D19').Value = '0.06944'
$ws.Range('E19').Value = '  +0.74%  '
$ws.Range('D20').Value = '17.56'
$ws.Range('E20').Value = '  -3.32%  '
$ws.Range('D21').Value = '6.531'
$ws.Range('E21').Value = '  +0.18%  '
$ws.Range('D22').Value = '0.9995'
$ws.Range('E22').Value = '  -0.06%  '
$ws.Range('E23').Value = '  -2.17%  '
$ws.Range('D24').Value = '23.383.77'
$ws.Range('E24').Value = '  -0.23%  '
$ws.Range('D25').Value = '2.511'
$ws.Range('E25').Value = '  +4.03%  '
$ws.Range('D26').Value = '3.077'
$ws.Range('E26').Value = '  -1.16%  '
$ws.Range('D27').Value = '21.14'
$ws.Range('E27').Value = '  -0.36%  '
$ws.Range('D28').Value = '150.45'
$ws.Range('E28').Value = '  -0.34%  '
$ws.Range('D29').Value = '5.269'
$ws.Range('E29').Value = '  -1.18%  '
$ws.Range('D30').Value = '132.86'
$ws.Range('E30').Value = '  -2.45%  '
$ws.Range('D31').Value = '1.800.08'
$ws.Range('E31').Value = '  -0.68%  '
$ws.Range('D32').Value = '6.617'
$ws.Range('E32').Value = '  -3.11%  '
$ws.Range('D33').Value = '2.137'
$ws.Range('E33').Value = '  -6.30%  '
$ws.Range('D34').Value = '1.058'
$ws.Range('E34').Value = '  +11.28%  '
$ws.Range('D35').Value = '11.34'
$ws.Range('E35').Value = '  +8.76%  '
$ws.Range('D36').Value = '0.02766'
$ws.Range('E36').Value = '  -1.65%  '
$ws.Range('E37').Value = '  -1.23%  '
$ws.Range('D38').Value = '0.08753'
$ws.Range('E38').Value = '  -0.08%  '
$ws.Range('D39').Value = '0.07135'
$ws.Range('E39').Value = '  -1.51%  '
$ws.Range('D40').Value = '5.978'
$ws.Range('E40').Value = '  -2.43%  '
$ws.Range('D41').Value = '0.6985'
$ws.Range('E41').Value = '  -1.26%  '
$ws.Range('D42').Value = '1.329'
$ws.Range('E42').Value = '  -3.01%  '
$ws.Range('D43').Value = '15.86'
$ws.Range('E43').Value = '  -1.36%  '
$ws.Range('D44').Value = '12.00'
$ws.Range('E44').Value = '  -4.09%  '
$ws.Range('D45').Value = '0.6463'
$ws.Range('E45').Value = '  -0.99%  '
$ws.Range('D46').Value = '0.9990'
$ws.Range('D47').Value = '2.275'
$ws.Range('E47').Value = '  -2.45%  '
$ws.Range('E48').Value = '  -1.41%  '
$ws.Range('D49').Value = '0.07974'
$ws.Range('E49').Value = '  +0.03%  '
$ws.Range('D50').Value = '126.41'
$ws.Range('E50').Value = '  -1.86%  '
$ws.Range('D51').Value = '1.187'
$ws.Range('E51').Value = '  -1.16%  '
